$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Haupt 2017): fill in notes/figures/organization/extracted fields
$ws.Range("M9").Value = "acclimation….maybe also included in the full analysis, table 1 excluded because there are no errors reported"
$ws.Range("K9").Value = "table 2, figure 1, 3a,  4"

# Row 7 (Bozinovic 2016): extend the "organization" note text with more detail
$ws.Range("L7").Value = "not sure how to handle constant mean with changing variance….so far, just extracted constant mean and constant variance trts as paired "

# Row 9 continued
$ws.Range("L9").Value = "not sure how to handle the data where they include and exclude individuals with different thermal preferences; decided to just use all individuals and not the exclusion analysis because I would be representing a lot of data twice and that might skew the data; extracted all data from figure 1 but not sure whether I should be including thermal preference data the both includes or excludes ctmin; assumed median values as means here "
$ws.Range("N9").Value = "y"

# Update the selected range shown in the sheet view
$ws.Range("N2:N9").Select()
